# ZoneAffix.xlsx update: EA 23.252 -> EA 23.260 (Chinese language file)
# - "unchanging"/"うつろわざる" affix (id 24) becomes "warped"/"歪んだ" and its
#   version bumps to EA 23.253; Chinese name becomes 扭曲的.
# - Ten new zone affixes (ids 31-40) are appended as new rows 33-42, all
#   stamped with version EA 23.253.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 26 (id 24): unchanging/うつろわざる -> warped/歪んだ ---
$ws.Cells.Item(26, 2).Value = "EA 23.253"
$ws.Cells.Item(26, 3).Value = "扭曲的"
$ws.Cells.Item(26, 4).Value = "warped"
$ws.Cells.Item(26, 5).Value = "歪んだ"

# --- Append the ten new affix rows (33-42), ids 31-40 ---
# Column A holds the numeric id as TEXT (matches existing id column formatting),
# so pre-format the range as text before writing the values.
$ws.Range("A33:A42").NumberFormat = "@"

$newRows = @(
    @{ Id = "31"; Ver = "EA 23.253"; CN = "亵渎的";     EN = "blasphemous"; JP = "冒涜的な" },
    @{ Id = "32"; Ver = "EA 23.253"; CN = "荒废的";     EN = "ruined";      JP = "荒廃した" },
    @{ Id = "33"; Ver = "EA 23.253"; CN = "崩塌的";     EN = "crumbling";   JP = "崩れかけた" },
    @{ Id = "34"; Ver = "EA 23.253"; CN = "与世隔绝的"; EN = "isolated";    JP = "隔絶された" },
    @{ Id = "35"; Ver = "EA 23.253"; CN = "充满魔力的"; EN = "arcane";      JP = "魔力に満ちた" },
    @{ Id = "36"; Ver = "EA 23.253"; CN = "生机勃勃的"; EN = "lively";      JP = "生き生きとした" },
    @{ Id = "37"; Ver = "EA 23.253"; CN = "朴素的";     EN = "rustic";      JP = "素朴な" },
    @{ Id = "38"; Ver = "EA 23.253"; CN = "引导的";     EN = "guiding";     JP = "導きの" },
    @{ Id = "39"; Ver = "EA 23.253"; CN = "泥土芬芳的"; EN = "earthy";      JP = "土の香る" },
    @{ Id = "40"; Ver = "EA 23.253"; CN = "受到庇护的"; EN = "warded";      JP = "護られた" }
)

$row = 33
foreach ($item in $newRows) {
    $ws.Cells.Item($row, 1).Value = $item.Id
    $ws.Cells.Item($row, 2).Value = $item.Ver
    $ws.Cells.Item($row, 3).Value = $item.CN
    $ws.Cells.Item($row, 4).Value = $item.EN
    $ws.Cells.Item($row, 5).Value = $item.JP
    $row++
}
